# RF010 - Gerenciar Periodos Avaliativos
# "From 0.1 to 1.0 version": bump the changelog entry and reword every
# remaining "perfis de competencias" reference to "Periodos Avaliativos"
# (plus a couple of related wording/footnote tweaks).

$d = $word.ActiveDocument

function Replace-FirstMatch {
    param(
        [string]$OldText,
        [string]$NewText
    )

    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $found) {
        throw "Could not find text: $OldText"
    }
    $rng.Text = $NewText
}

# 1) Changelog table: version bump 0.1 -> 1.0
Replace-FirstMatch "0.1" "1.0"

# 2) Changelog table: change description Criacao -> Final
Replace-FirstMatch "Criacao" "Final"

# 3) Precondition cell (text runs on either side of an embedded line
#    break, so the search term is kept inside a single "line" to avoid
#    Find crossing the literal CR/LF characters)
Replace-FirstMatch `
    "tem permissao para gerenciar perfis de competencias." `
    "tem permissao para gerenciar Periodos Avaliativos."

# 4) Main flow step 1
Replace-FirstMatch `
    "1. Lider de Pessoas acessa a funcionalidade de gestao de perfis de competencias a partir do menu inicial af[4]" `
    "1. Lider de Pessoas acessa a funcionalidade de gestao de Periodos Avaliativos a partir do menu inicial af[4]"

# 5) Main flow step 2
Replace-FirstMatch `
    "2. System exibe a listagem dos perfis de competencias cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda' " `
    "2. System exibe a listagem dos Periodos Avaliativos cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda' "

# 6) Main flow step 7
Replace-FirstMatch `
    "7. Lider de Pessoas preenche o campo 'Data Inicial e Data Final' selecionando um lider da lista " `
    "7. Lider de Pessoas preenche o campo 'Data Inicial e Data Final' informando as respectivas datas referente ao periodo "

# 7) Main flow step 18 footnote reference (also separated from the
#    embedded line break the same way as change 3)
Replace-FirstMatch `
    "retorna feedback correspondente ef[1,2,3]" `
    "retorna feedback correspondente ef[1,2]"

# 8) EF[3] flow step 6 (also gains a footnote reference)
Replace-FirstMatch `
    "6. System exibe a listagem dos perfis de competencias sem o Periodo Avaliativo excluido " `
    "6. System exibe a listagem dos Periodos Avaliativos sem o Periodo Avaliativo excluido ef[3]"

# 9) AF flow step 6
Replace-FirstMatch `
    "6. System exibe a listagem dos perfis de competencias com o Periodo Avaliativo excluido " `
    "6. System exibe a listagem dos Periodos Avaliativos com o Periodo Avaliativo excluido "

# 10) Unauthenticated-user flow step 1
Replace-FirstMatch `
    "1. Usuario Nao-Autenticado acessa a funcionalidade de gestao de perfis de competencias a partir do menu inicial " `
    "1. Usuario Nao-Autenticado acessa a funcionalidade de gestao de Periodos Avaliativos a partir do menu inicial "

# 11) Unauthenticated-user flow step 2
Replace-FirstMatch `
    "2. System exibe a listagem dos perfis de competencias cadastrados apenas para visualizacao com a opcao 'Ajuda' " `
    "2. System exibe a listagem dos Periodos Avaliativos cadastrados apenas para visualizacao com a opcao 'Ajuda' "

# 12) Post-conditions table
Replace-FirstMatch `
    "A gestao de perfis de competencias e realizada com sucesso." `
    "A gestao de Periodos Avaliativos e realizada com sucesso."
